$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (TRIPLIXAM 10/2.5/10MG 15 F.C. TABS.)
# Stock text changes from "-90:0" to "0:0"
$ws.Range("H9").Value = "0:0"

# Sell price text changes from "21294.0000" to "234.0000"
# (value looks numeric, so force text so it stays a shared string like the original)
$fmtP9 = $ws.Range("P9").NumberFormat
$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "234.0000"
$ws.Range("P9").NumberFormat = $fmtP9

# Transaction count text changes from "91:0" to "1:0"
$ws.Range("Q9").Value = "1:0"

# Updated total in P12
$ws.Range("P12").Value = 333.83

# Updated generation timestamp
$ws.Range("A13").Value = "Wednesday, 3 September, 2025 9:57 AM"
